$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => (Median Value, Tier) after recalculating scores relative to the median
# AFTER merging with zip/census tract data (see commit message).
$rows = @(
    @(2, 0.5899937067337949, "Below Median"),
    @(3, 0.3477029578351165, "Below Median"),
    @(4, 1.75330396475771, "1st Tier"),
    @(5, 1.874082232011748, "1st Tier"),
    @(6, 2.314610866372981, "1st Tier"),
    @(7, 1.198867212083071, "3rd Tier"),
    @(8, 1.037444933920705, "4th Tier"),
    @(9, 2.279735682819383, "1st Tier"),
    @(10, 0.5711139081183135, "Below Median"),
    @(11, 0.8039647577092511, "Below Median"),
    @(12, 0.9770295783511644, "Below Median"),
    @(13, 1.500550660792952, "2nd Tier"),
    @(14, 0.7544052863436124, "Below Median"),
    @(15, 1.082966226138032, "4th Tier"),
    @(16, 1, "4th Tier"),
    @(17, 0.9512665198237886, "Below Median"),
    @(18, 0.5451541850220265, "Below Median"),
    @(19, 0.938570729319628, "Below Median"),
    @(20, 0.7672540381791484, "Below Median"),
    @(21, 0.3634361233480177, "Below Median"),
    @(22, 0.8731906859660165, "Below Median"),
    @(23, 0.8763373190685967, "Below Median"),
    @(24, 1.651982378854626, "1st Tier"),
    @(25, 1.342039018250472, "3rd Tier"),
    @(26, 1.5712187958884, "2nd Tier"),
    @(27, 1.506057268722467, "2nd Tier"),
    @(28, 1.120594713656388, "4th Tier"),
    @(29, 0.7158590308370044, "Below Median"),
    @(30, 1.82488986784141, "1st Tier"),
    @(31, 1.260401370533529, "3rd Tier"),
    @(32, 1.583700440528634, "2nd Tier"),
    @(33, 0.4937591776798826, "Below Median"),
    @(34, 1.135934550031467, "4th Tier"),
    @(35, 0.6718061674008811, "Below Median"),
    @(36, 1.551290119572058, "2nd Tier"),
    @(37, 0.9990560100692261, "Below Median"),
    @(38, 1.593769666456891, "2nd Tier"),
    @(39, 1.173694147262429, "3rd Tier"),
    @(40, 1.814977973568282, "1st Tier"),
    @(41, 1.307819383259912, "3rd Tier"),
    @(42, 1.042584434654919, "4th Tier"),
    @(43, 1.668502202643172, "1st Tier"),
    @(44, 0.9040849018822588, "Below Median"),
    @(45, 0.8513215859030838, "Below Median"),
    @(46, 0.1916299559471366, "Below Median"),
    @(47, 0.947136563876652, "Below Median"),
    @(48, 1.049008810572687, "4th Tier"),
    @(49, 0.9295154185022027, "Below Median"),
    @(50, 0.9691629955947137, "Below Median"),
    @(51, 0.5855359765051396, "Below Median"),
    @(52, 1.02863436123348, "4th Tier"),
    @(53, 1.613436123348018, "2nd Tier"),
    @(54, 0.7075991189427313, "Below Median"),
    @(55, 1.458464443045941, "2nd Tier"),
    @(56, 1.202013845185651, "3rd Tier"),
    @(57, 0.4118942731277533, "Below Median"),
    @(58, 0.7632158590308371, "Below Median"),
    @(59, 0.3193832599118943, "Below Median"),
    @(60, 0.6104468219005664, "Below Median"),
    @(61, 1.236233480176212, "3rd Tier"),
    @(62, 1.698788546255507, "1st Tier"),
    @(63, 0.8072687224669605, "Below Median"),
    @(64, 0.6916299559471366, "Below Median"),
    @(65, 0.9361233480176212, "Below Median"),
    @(66, 0.4717327459618209, "Below Median"),
    @(67, 1.411894273127754, "2nd Tier"),
    @(68, 1.170154185022027, "3rd Tier"),
    @(69, 0.7246696035242292, "Below Median"),
    @(70, 0.7709251101321587, "Below Median"),
    @(71, 1.079295154185022, "4th Tier"),
    @(72, 1.288546255506608, "3rd Tier"),
    @(73, 0.7378854625550662, "Below Median"),
    @(74, 1.587476400251731, "2nd Tier"),
    @(75, 2.702643171806168, "1st Tier"),
    @(76, 0.8241556534508077, "Below Median"),
    @(77, 1.110132158590308, "4th Tier"),
    @(78, 1.183920704845815, "3rd Tier"),
    @(79, 0.6084801762114538, "Below Median"),
    @(80, 0.7929515418502203, "Below Median")
)

foreach ($row in $rows) {
    $r = $row[0]
    $medianValue = $row[1]
    $tier = $row[2]
    $ws.Cells.Item($r, 3).Value = $medianValue
    $ws.Cells.Item($r, 4).Value = $tier
}
